# Update computed profit/price figures across sheets per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H40").Value = 8896.666999999999
$ws_ALC.Range("I40").Value = 12900
$ws_ALC.Range("K40").Value = 12900
$ws_ALC.Range("M40").Value = -12725
$ws_ALC.Range("H133").Value = 14999.667
$ws_ALC.Range("J133").Value = 14999.667
$ws_ALC.Range("L133").Value = 14999.667
$ws_ALC.Range("N133").Value = -25119.667
$ws_ALC.Range("H135").Value = 4052.4546
$ws_ALC.Range("I135").Value = 4207.7
$ws_ALC.Range("J135").Value = 2500
$ws_ALC.Range("K135").Value = 37869.3
$ws_ALC.Range("L135").Value = 22500
$ws_ALC.Range("M135").Value = -35334.3
$ws_ALC.Range("N135").Value = -27570
$ws_ALC.Range("H136").Value = 67673.62
$ws_ALC.Range("J136").Value = 67673.62
$ws_ALC.Range("L136").Value = 67673.62
$ws_ALC.Range("N136").Value = -77873.62
$ws_ALC.Range("H138").Value = 2139.77
$ws_ALC.Range("I138").Value = 1356.5
$ws_ALC.Range("J138").Value = 2360.6924
$ws_ALC.Range("K138").Value = 4069.5
$ws_ALC.Range("L138").Value = 7082.0772
$ws_ALC.Range("M138").Value = 1070.5
$ws_ALC.Range("N138").Value = -17362.0772

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H34").Value = 35000
$ws_ARM.Range("J34").Value = 35000
$ws_ARM.Range("L34").Value = 35000
$ws_ARM.Range("N34").Value = -35542
$ws_ARM.Range("H45").Value = 3714.8928
$ws_ARM.Range("I45").Value = 3689.4736
$ws_ARM.Range("K45").Value = 3689.4736
$ws_ARM.Range("M45").Value = -3312.4736
$ws_ARM.Range("H132").Value = 20713.115
$ws_ARM.Range("I132").Value = 5825.5
$ws_ARM.Range("J132").Value = 33473.93
$ws_ARM.Range("K132").Value = 17476.5
$ws_ARM.Range("L132").Value = 100421.79
$ws_ARM.Range("M132").Value = -14946.5
$ws_ARM.Range("N132").Value = -105481.79
$ws_ARM.Range("H140").Value = 88333.336
$ws_ARM.Range("J140").Value = 88333.336
$ws_ARM.Range("L140").Value = 88333.336
$ws_ARM.Range("N140").Value = -98693.336

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H123").Value = 0
$ws_BSM.Range("J123").Value = 0
$ws_BSM.Range("L123").Value = 0
$ws_BSM.Range("N123").ClearContents()
$ws_BSM.Range("H134").Value = 44888.242
$ws_BSM.Range("I134").Value = 48537.773
$ws_BSM.Range("J134").Value = 33418.285
$ws_BSM.Range("K134").Value = 145613.319
$ws_BSM.Range("L134").Value = 100254.855
$ws_BSM.Range("M134").Value = -143078.319
$ws_BSM.Range("N134").Value = -105324.855
$ws_BSM.Range("H137").Value = 83000
$ws_BSM.Range("J137").Value = 83000
$ws_BSM.Range("L137").Value = 83000
$ws_BSM.Range("N137").Value = -93200
$ws_BSM.Range("H140").Value = 78398.7
$ws_BSM.Range("J140").Value = 78398.7
$ws_BSM.Range("L140").Value = 78398.7
$ws_BSM.Range("N140").Value = -88758.7

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H23").Value = 7366.3335
$ws_CRP.Range("I23").Value = 7499
$ws_CRP.Range("J23").Value = 7300
$ws_CRP.Range("K23").Value = 7499
$ws_CRP.Range("L23").Value = 7300
$ws_CRP.Range("M23").Value = -7259
$ws_CRP.Range("N23").Value = -7780
$ws_CRP.Range("H27").Value = 7366.3335
$ws_CRP.Range("I27").Value = 7499
$ws_CRP.Range("J27").Value = 7300
$ws_CRP.Range("K27").Value = 7499
$ws_CRP.Range("L27").Value = 7300
$ws_CRP.Range("M27").Value = -7307
$ws_CRP.Range("N27").Value = -7684
$ws_CRP.Range("H74").Value = 35000
$ws_CRP.Range("J74").Value = 35000
$ws_CRP.Range("L74").Value = 35000
$ws_CRP.Range("N74").Value = -36748
$ws_CRP.Range("H77").Value = 35000
$ws_CRP.Range("J77").Value = 35000
$ws_CRP.Range("L77").Value = 105000
$ws_CRP.Range("N77").Value = -113736
$ws_CRP.Range("H134").Value = 8852.299999999999
$ws_CRP.Range("I134").Value = 2584.9167
$ws_CRP.Range("J134").Value = 18253.375
$ws_CRP.Range("K134").Value = 7754.750100000001
$ws_CRP.Range("L134").Value = 54760.125
$ws_CRP.Range("M134").Value = -5219.750100000001
$ws_CRP.Range("N134").Value = -59830.125

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H34").Value = 1409.7
$ws_CUL.Range("I34").Value = 510.77777
$ws_CUL.Range("J34").Value = 9500
$ws_CUL.Range("K34").Value = 1532.33331
$ws_CUL.Range("L34").Value = 28500
$ws_CUL.Range("M34").Value = -1448.33331
$ws_CUL.Range("N34").Value = -28668
$ws_CUL.Range("H38").Value = 106.75
$ws_CUL.Range("I38").Value = 52.25
$ws_CUL.Range("J38").Value = 215.75
$ws_CUL.Range("K38").Value = 156.75
$ws_CUL.Range("L38").Value = 647.25
$ws_CUL.Range("M38").Value = 190.25
$ws_CUL.Range("N38").Value = -1341.25
$ws_CUL.Range("H39").Value = 8283.166999999999
$ws_CUL.Range("I39").Value = 5099.5
$ws_CUL.Range("J39").Value = 9875
$ws_CUL.Range("K39").Value = 15298.5
$ws_CUL.Range("L39").Value = 29625
$ws_CUL.Range("M39").Value = -15004.5
$ws_CUL.Range("N39").Value = -30213
$ws_CUL.Range("H55").Value = 6211.2856
$ws_CUL.Range("I55").Value = 5095.8
$ws_CUL.Range("J55").Value = 9000
$ws_CUL.Range("K55").Value = 15287.4
$ws_CUL.Range("L55").Value = 27000
$ws_CUL.Range("M55").Value = -15110.4
$ws_CUL.Range("N55").Value = -27354
$ws_CUL.Range("H92").Value = 408.63635
$ws_CUL.Range("J92").Value = 459
$ws_CUL.Range("L92").Value = 1377
$ws_CUL.Range("N92").Value = -3873
$ws_CUL.Range("H113").Value = 996.8929000000001
$ws_CUL.Range("J113").Value = 1038.3125
$ws_CUL.Range("L113").Value = 3114.9375
$ws_CUL.Range("N113").Value = -7454.9375

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 14166.546
$ws_GSM.Range("I80").Value = 13671.083
$ws_GSM.Range("K80").Value = 13671.083
$ws_GSM.Range("M80").Value = -12673.083
$ws_GSM.Range("H83").Value = 14166.546
$ws_GSM.Range("I83").Value = 13671.083
$ws_GSM.Range("K83").Value = 68355.41500000001
$ws_GSM.Range("M83").Value = -63363.41500000001
$ws_GSM.Range("H103").Value = 87750
$ws_GSM.Range("J103").Value = 87750
$ws_GSM.Range("L103").Value = 87750
$ws_GSM.Range("N103").Value = -90094
$ws_GSM.Range("H132").Value = 27189.5
$ws_GSM.Range("I132").Value = 19097.889
$ws_GSM.Range("J132").Value = 100014
$ws_GSM.Range("K132").Value = 57293.667
$ws_GSM.Range("L132").Value = 300042
$ws_GSM.Range("M132").Value = -54763.667
$ws_GSM.Range("N132").Value = -305102
$ws_GSM.Range("H134").Value = 65636.2
$ws_GSM.Range("J134").Value = 65636.2
$ws_GSM.Range("L134").Value = 196908.6
$ws_GSM.Range("N134").Value = -201978.6

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H4").Value = 14666.667
$ws_LTW.Range("I4").Value = 14000
$ws_LTW.Range("K4").Value = 14000
$ws_LTW.Range("M4").Value = -13887
$ws_LTW.Range("H28").Value = 14666.667
$ws_LTW.Range("I28").Value = 14000
$ws_LTW.Range("K28").Value = 14000
$ws_LTW.Range("M28").Value = -13768
$ws_LTW.Range("H37").Value = 14666.667
$ws_LTW.Range("I37").Value = 14000
$ws_LTW.Range("K37").Value = 14000
$ws_LTW.Range("M37").Value = -13893
$ws_LTW.Range("H68").Value = 3432.8
$ws_LTW.Range("I68").Value = 3522
$ws_LTW.Range("K68").Value = 3522
$ws_LTW.Range("M68").Value = -2773
$ws_LTW.Range("H71").Value = 3432.8
$ws_LTW.Range("I71").Value = 3522
$ws_LTW.Range("K71").Value = 17610
$ws_LTW.Range("M71").Value = -13866
$ws_LTW.Range("H96").Value = 17000
$ws_LTW.Range("J96").Value = 17000
$ws_LTW.Range("L96").Value = 17000
$ws_LTW.Range("N96").Value = -22492
$ws_LTW.Range("H132").Value = 19194.584
$ws_LTW.Range("J132").Value = 24632.111
$ws_LTW.Range("L132").Value = 73896.333
$ws_LTW.Range("N132").Value = -78956.333

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H121").Value = 0
$ws_WVR.Range("J121").Value = 0
$ws_WVR.Range("L121").Value = 0
$ws_WVR.Range("N121").ClearContents()
$ws_WVR.Range("H124").Value = 32379.8
$ws_WVR.Range("J124").Value = 32379.8
$ws_WVR.Range("L124").Value = 32379.8
$ws_WVR.Range("N124").Value = -42199.8
$ws_WVR.Range("H132").Value = 6660.8335
$ws_WVR.Range("J132").Value = 14998.615
$ws_WVR.Range("L132").Value = 44995.845
$ws_WVR.Range("N132").Value = -50055.845
$ws_WVR.Range("H136").Value = 15494.6
$ws_WVR.Range("I136").Value = 1815.2142
$ws_WVR.Range("J136").Value = 32904.727
$ws_WVR.Range("K136").Value = 5445.642599999999
$ws_WVR.Range("L136").Value = 98714.181
$ws_WVR.Range("M136").Value = -2895.642599999999
$ws_WVR.Range("N136").Value = -103814.181
